# Update "想去人数" (number of interested people) values for a few events
# across the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 55
$wsExhibit.Range("F5").Value = 5032
$wsExhibit.Range("F9").Value = 51

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 55
$wsAll.Range("F9").Value = 5032
$wsAll.Range("F14").Value = 51
